# New5VIN_UT_SS.xlsx — "added individual vins for the ss tests: pt3 SS"
#
# The shared VIN value used by rows 2-5 (column A) changes from
# "TTTKN3DD&E" to "EEEKN3DD&E". All four cells share the same string,
# so all four must be rewritten with the new value for the workbook's
# shared-string table to collapse back down to a single entry (matching
# the one-line diff against sharedStrings.xml).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "EEEKN3DD&E"
$ws.Range("A3").Value = "EEEKN3DD&E"
$ws.Range("A4").Value = "EEEKN3DD&E"
$ws.Range("A5").Value = "EEEKN3DD&E"

# The author's last selection before saving moved from C11 to B11.
[void]$ws.Range("B11").Select()
